# Update lzs vocab namespace en conceptscheme base uri
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header typo + simplify the referentiewaardeType header value.
$ws.Range("W1").Value = "technischefiche"
$ws.Range("X1").Value = "referentiewaardeType"

# All conceptscheme base URIs (columns A, H, J) need an extra "/lzs"
# path segment inserted right after ".../id/conceptscheme/", e.g.
#   https://data.omgeving.vlaanderen.be/id/conceptscheme/lzsp
# becomes
#   https://data.omgeving.vlaanderen.be/id/conceptscheme/lzs/lzsp
$oldBase = "/id/conceptscheme/"
$newBase = "/id/conceptscheme/lzs/"

$lastRow = $ws.UsedRange.Rows.Count
$cols = @(1, 8, 10)   # A = 1 (concept/scheme id), H = 8 (inScheme), J = 10 (topConceptOf)

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($c in $cols) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -ne $null -and $val -like "*$oldBase*" -and $val -notlike "*$newBase*") {
            $cell.Value = $val -replace [regex]::Escape($oldBase), $newBase
        }
    }
}
